# Export with no is_pref and no lev distance
#
# The "id" column (B) used to be produced by some fuzzy / Levenshtein-distance
# matching step and the "is_prefered" column (D) flagged the chosen match with
# an "x". The new export simply derives the id directly from the
# speaker_variant (column C) as "#" + lowercase(speaker_variant), and no
# longer marks any row as preferred, so column D is left blank for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final speaker_variant (column C) values for each data row, in row order,
# as produced by the new export (rows were also re-ordered/re-grouped by the
# updated matching logic upstream).
$variants = @{
    2  = "Stra"
    3  = "Lyc"
    4  = "Arch"
    5  = "Pol"
    6  = "App"
    7  = "Sich"
    8  = "Eld"
    9  = "Dard"
    10 = "Antio"
    11 = "Pygm"
    12 = "Ant"
    13 = "Barg"
    14 = "Lye"
    15 = "Arg"
    16 = "Had"
    17 = "Strag"
    18 = "Krat"
    19 = "Elis"
    20 = "Aret"
    21 = "Ptol"
    22 = "Ar"
}

foreach ($row in 2..22) {
    $variant = $variants[$row]

    # id = "#" + lowercase(speaker_variant) -- no more lev-distance based id
    $ws.Cells.Item($row, 2).Value = "#" + $variant.ToLower()

    # speaker_variant
    $ws.Cells.Item($row, 3).Value = $variant

    # is_prefered -- no longer exported, clear any previous "x" marker
    $ws.Cells.Item($row, 4).ClearContents()
}
